# Update macrophage annotation labels
# - Rename header "Broad" -> "ann_level_1", "Annotation" -> "ann_level_2"
# - Insert a new "ann_level_3" column (D) holding the original, detailed
#   cluster annotation that used to live in column C
# - Column C becomes the broader cluster-family label (e.g. "macro-CCL18"
#   collapses to "macro-CCL", "macro-proliferating-S" collapses to
#   "macro-proliferating", etc.)
# - Column B ("Broad") values are lower-cased ("Macrophages" -> "macrophages",
#   "Proliferating macrophages" -> "proliferating macrophages")
# - "Relevant marker genes"/"Relevant marker ADTs" columns shift right one
#   column to make room (old D/E/F -> E/F/G)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at D; this shifts the old D/E/F (marker genes,
# marker ADTs, reactome) columns to E/F/G and carries their widths/styles.
$ws.Columns.Item(4).Insert()

# ---- Header row ----
$ws.Cells.Item(1, 2).Value = "ann_level_1"
$ws.Cells.Item(1, 3).Value = "ann_level_2"
$ws.Cells.Item(1, 4).Value = "ann_level_3"

# ---- Data rows ----
# Each entry: row, ann_level_1 (B), ann_level_2 (C, broad cluster), ann_level_3 (D, detailed cluster)
$rows = @(
    @(2,  "macrophages",               "macro-alveolar",         "macro-alveolar"),
    @(3,  "macrophages",               "macro-IFI27",            "macro-IFI27"),
    @(4,  "macrophages",               "macro-monocyte-derived", "macro-monocyte-derived"),
    @(5,  "macrophages",               "macro-APOC2+",           "macro-APOC2+"),
    @(6,  "macrophages",               "macro-alveolar",         "macro-alveolar"),
    @(7,  "macrophages",               "macro-alveolar",         "macro-alveolar"),
    @(8,  "macrophages",               "macro-CCL",              "macro-CCL18"),
    @(9,  "macrophages",               "macro-CCL",              "macro-CCL"),
    @(10, "macrophages",               "macro-lipid",             "macro-lipid"),
    @(11, "macrophages",               "macro-IGF1",              "macro-IGF1"),
    @(12, "macrophages",               "macro-IFI27",             "macro-IFI27+APOC2+"),
    @(13, "proliferating macrophages", "macro-proliferating",     "macro-proliferating-S"),
    @(14, "macrophages",               "macro-MT",                "macro-MT"),
    @(15, "macrophages",               "macro-interstitial",      "macro-interstitial"),
    @(16, "macrophages",               "macro-T",                 "macro-T"),
    @(17, "proliferating macrophages", "macro-proliferating",     "macro-proliferating-G2M"),
    @(18, "macrophages",               "macro-lipid",              "macro-lipid-APOC2+"),
    @(19, "macrophages",               "macro-IFI27",              "macro-IFI27+CCL18+"),
    @(20, "macrophages",               "macro-IFN",                "macro-IFN"),
    @(21, "unknown",                   "unknown",                  "unknown"),
    @(22, "unknown",                   "unknown",                  "unknown"),
    @(23, "proliferating macrophages", "macro-proliferating",      "macro-proliferating-G2M")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
}

# Match the author's last selection after the edit.
$ws.Range("D26").Select() | Out-Null

Write-Output "macrophage annotation labels updated"
